$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "35.416.91"
Set-TextValue "E2" "  +0.33%  "
Set-TextValue "D3" "1.892.99"
Set-TextValue "E3" "  -0.98%  "
Set-TextValue "E4" "  -0.83%  "
Set-TextValue "D5" "247.37"
Set-TextValue "E5" "  -2.54%  "
Set-TextValue "D6" "0.689"
Set-TextValue "E6" "  -4.70%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.90%  "
Set-TextValue "D8" "44.00"
Set-TextValue "E8" "  +8.17%  "
Set-TextValue "D9" "0.353"
Set-TextValue "E9" "  -2.61%  "
Set-TextValue "D10" "51.92"
Set-TextValue "E10" "  -0.65%  "
Set-TextValue "D11" "0.0739"
Set-TextValue "E11" "  -3.87%  "
Set-TextValue "E12" "  -1.68%  "
Set-TextValue "D13" "13.20"
Set-TextValue "E13" "  +3.31%  "
Set-TextValue "D14" "2.161.78"
Set-TextValue "E14" "  -1.28%  "
Set-TextValue "E15" "  +0.45%  "
Set-TextValue "D16" "1.914.19"
Set-TextValue "E16" "  +1.26%  "
Set-TextValue "D17" "4.93"
Set-TextValue "E17" "  -0.13%  "
Set-TextValue "D18" "35.384.78"
Set-TextValue "E18" "  +0.22%  "
Set-TextValue "D19" "73.22"
Set-TextValue "E19" "  -1.61%  "
Set-TextValue "D20" "0.0₃0822"
Set-TextValue "E20" "  -3.32%  "
Set-TextValue "D21" "246.05"
Set-TextValue "E21" "  +0.84%  "
Set-TextValue "D22" "12.85"
Set-TextValue "E22" "  -1.86%  "
Set-TextValue "E23" "  -2.22%  "
Set-TextValue "E24" "  -0.90%  "
Set-TextValue "D25" "2.55"
Set-TextValue "E25" "  +6.83%  "
Set-TextValue "D26" "2.20"
Set-TextValue "E26" "  -10.50%  "
Set-TextValue "D27" "165.56"
Set-TextValue "E27" "  -0.77%  "
Set-TextValue "D28" "8.51"
Set-TextValue "E28" "  -1.89%  "
Set-TextValue "D29" "18.32"
Set-TextValue "E29" "  -2.23%  "
Set-TextValue "E30" "  -4.21%  "
Set-TextValue "D31" "4.128.46"
Set-TextValue "E31" "  +0.04%  "
Set-TextValue "D32" "1.79"
Set-TextValue "E32" "  +9.64%  "
Set-TextValue "D33" "4.28"
Set-TextValue "E33" "  -1.51%  "
Set-TextValue "D34" "0.0581"
Set-TextValue "E34" "  -0.44%  "
Set-TextValue "D35" "4.25"
Set-TextValue "E35" "  +0.68%  "
Set-TextValue "E36" "  -0.88%  "
Set-TextValue "D37" "0.849"
Set-TextValue "E37" "  -7.45%  "
Set-TextValue "E38" "  -1.71%  "
Set-TextValue "D39" "1.59"
Set-TextValue "E39" "  -20.43%  "
Set-TextValue "D40" "17.25"
Set-TextValue "E40" "  +0.41%  "
Set-TextValue "D41" "97.62"
Set-TextValue "E41" "  +1.01%  "
Set-TextValue "D42" "0.0670"
Set-TextValue "D43" "0.0213"
Set-TextValue "E43" "  -1.95%  "
Set-TextValue "E44" "  -1.88%  "
Set-TextValue "D45" "1.292.06"
Set-TextValue "E45" "  -3.49%  "
Set-TextValue "D46" "2.37"
Set-TextValue "E46" "  -2.33%  "
Set-TextValue "E47" "  +7.98%  "
Set-TextValue "E48" "  -0.91%  "
Set-TextValue "D49" "2.75"
Set-TextValue "E49" "  -0.88%  "
Set-TextValue "D50" "12.05"
Set-TextValue "E50" "  +0.23%  "
Set-TextValue "D51" "6.41"
Set-TextValue "E51" "  -4.95%  "
